$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value2 = $statusText
$overview.Range("C2").Value2 = $statusText
$overview.Range("B3").Value2 = $statusText
$overview.Range("C3").Value2 = $statusText

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value2 = $statusText
$zhcn.Range("C3").Value2 = $statusText
$zhcn.Range("H2").Value2 = "2016-03-22 15:09:20"
$zhcn.Range("H3").Value2 = "2016-03-22 15:09:20"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value2 = $statusText
$dede.Range("C3").Value2 = $statusText
$dede.Range("H2").Value2 = "2016-03-22 15:09:26"
$dede.Range("H3").Value2 = "2016-03-22 15:09:26"

Write-Output "status+datetime updates done"

function Add-HandbackColumns($ws, $aUrl, $dUrl, $a3Url, $d3Url) {
  $ws.Range("F2").Value2 = $ws.Range("A2").Value2
  $ws.Hyperlinks.Add($ws.Range("F2"), $aUrl, "", "", $ws.Range("A2").Value2) | Out-Null
  $ws.Range("F2").Style = "HyperLink"

  $ws.Range("G2").Value2 = $ws.Range("D2").Value2
  $ws.Hyperlinks.Add($ws.Range("G2"), $dUrl, "", "", $ws.Range("D2").Value2) | Out-Null
  $ws.Range("G2").Style = "HyperLink"

  $ws.Range("F3").Value2 = $ws.Range("A3").Value2
  $ws.Hyperlinks.Add($ws.Range("F3"), $a3Url, "", "", $ws.Range("A3").Value2) | Out-Null
  $ws.Range("F3").Style = "HyperLink"

  $ws.Range("G3").Value2 = $ws.Range("D3").Value2
  $ws.Hyperlinks.Add($ws.Range("G3"), $d3Url, "", "", $ws.Range("D3").Value2) | Out-Null
  $ws.Range("G3").Style = "HyperLink"
}

$mdUrlA = "https://github.com/OpenLocalizationTest/oltest/blob/cf6324a76420f0f0018fbd7dd1dea68225614eb5/e2e/a8a25305-d45e-4370-9a91-f732f20bbebb.md"
$mdUrlA3 = "https://github.com/OpenLocalizationTest/oltest/blob/cf6324a76420f0f0018fbd7dd1dea68225614eb5/e2e/cd7831bf-1df1-47ac-b351-225d4c84c32e.md"
$zhcnXlfUrlD = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20c2766b56366a64264ed2e9cee1ee34ede03256/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a8a25305-d45e-4370-9a91-f732f20bbebb.614e87983797ec001f34047d572375cf0a25c393.zh-cn.xlf"
$zhcnXlfUrlD3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20c2766b56366a64264ed2e9cee1ee34ede03256/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/cd7831bf-1df1-47ac-b351-225d4c84c32e.628f699113f9846d5c902626a8e900646dbf9d5e.zh-cn.xlf"
$dedeXlfUrlD = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9e624e4e13e357d6aeae30e77fc61140809d5a2b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a8a25305-d45e-4370-9a91-f732f20bbebb.614e87983797ec001f34047d572375cf0a25c393.de-de.xlf"
$dedeXlfUrlD3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9e624e4e13e357d6aeae30e77fc61140809d5a2b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/cd7831bf-1df1-47ac-b351-225d4c84c32e.628f699113f9846d5c902626a8e900646dbf9d5e.de-de.xlf"

Add-HandbackColumns $zhcn $mdUrlA $zhcnXlfUrlD $mdUrlA3 $zhcnXlfUrlD3
Add-HandbackColumns $dede $mdUrlA $dedeXlfUrlD $mdUrlA3 $dedeXlfUrlD3

Write-Output "handback columns done"
